$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B16: "None" -> "nan"
$ws.Range("B16").Value = "nan"

# New row 17. A17/D17 hold the numeric-looking text "123"; prefix with an
# apostrophe (the standard GUI way of forcing a number-like entry to stay
# text) and then restore the "Normal" style so no stray number-format/
# quote-prefix formatting is left behind on the cell.
$ws.Range("A17").Value = "'123"
$ws.Range("A17").Style = "Normal"

$ws.Range("B17").Value = "wew"

$ws.Range("C17").Value = "sd"

$ws.Range("D17").Value = "'123"
$ws.Range("D17").Style = "Normal"
